$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 240, shifting existing rows 240-290 down to 241-291.
$ws.Rows(240).Insert()

# Populate the newly inserted row 240 with the new record's data.
$ws.Range("A240").Value = 4
$ws.Range("B240").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C240").Value = "Los Lagos"
$ws.Range("D240").Value = 44782
$ws.Range("E240").Value = 10
$ws.Range("F240").Value = 100112021
$ws.Range("G240").Value = "Ají"
$ws.Range("H240").Value = "Inferno"
$ws.Range("I240").Value = "Primera"
$ws.Range("J240").Value = 140
$ws.Range("K240").Value = 21000
$ws.Range("L240").Value = 21000
$ws.Range("M240").Value = 21000
$ws.Range("N240").Value = "$/caja 12 kilos"
$ws.Range("O240").Value = "Región de Arica y Parinacota"
$ws.Range("P240").Value = 1750
$ws.Range("Q240").Value = 12
$ws.Range("R240").Value = "Hortaliza"
